# Update BunkerPrices at 2025-04-04 13:20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the "Hong Kong" / "Montreal*" columns (AB <-> AC) ---
# Header row (row 1)
$tmp = $ws.Range("AB1").Value()
$ws.Range("AB1").Value = $ws.Range("AC1").Value()
$ws.Range("AC1").Value = $tmp

# Data rows 2-16
for ($r = 2; $r -le 16; $r++) {
    $abRef = "AB" + $r
    $acRef = "AC" + $r
    $tmp = $ws.Range($abRef).Value()
    $ws.Range($abRef).Value = $ws.Range($acRef).Value()
    $ws.Range($acRef).Value = $tmp
}

# --- Row 16's date cell reverts to the normal date+time format ---
# (it was previously the last row and used the "latest row" date-only format)
$ws.Range("E16").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# --- Append the new data row 17 ---
$ws.Range("A17").Value = 565
$ws.Range("B17").Value = 585
$ws.Range("C17").Value = 878
$ws.Range("D17").Value = 532
$ws.Range("E17").Value = 45749
$ws.Range("F17").Value = 534
$ws.Range("G17").Value = 535
$ws.Range("H17").Value = 574
$ws.Range("I17").Value = 584
$ws.Range("J17").Value = 552
$ws.Range("K17").Value = 572
$ws.Range("L17").Value = 595
$ws.Range("M17").Value = 502
$ws.Range("N17").Value = 600
$ws.Range("O17").Value = 656
$ws.Range("P17").Value = 535
$ws.Range("Q17").Value = 584
$ws.Range("R17").Value = 532
$ws.Range("S17").Value = 564
$ws.Range("T17").Value = 668
$ws.Range("U17").Value = 623
$ws.Range("V17").Value = 603
$ws.Range("W17").Value = 625
$ws.Range("X17").Value = 534
$ws.Range("Y17").Value = 577
$ws.Range("Z17").Value = 778
$ws.Range("AA17").Value = 535
$ws.Range("AB17").Value = 675
$ws.Range("AC17").Value = 540
$ws.Range("AD17").Value = 629
$ws.Range("AE17").Value = 532
$ws.Range("AF17").Value = 565
$ws.Range("AG17").Value = 515
$ws.Range("AH17").Value = 667
$ws.Range("AI17").Value = 535
$ws.Range("AJ17").Value = 559
$ws.Range("AK17").Value = 619.5
$ws.Range("AL17").Value = 645
$ws.Range("AM17").Value = 545
$ws.Range("AN17").Value = 537
$ws.Range("AO17").Value = 580
$ws.Range("AP17").Value = 663
$ws.Range("AQ17").Value = 525
$ws.Range("AR17").Value = 556
$ws.Range("AS17").Value = 502
$ws.Range("AT17").Value = 542
$ws.Range("AU17").Value = 788
$ws.Range("AV17").Value = 530

# New last row gets the "latest row" date-only format
$ws.Range("E17").NumberFormat = "YYYY-MM-DD"
